$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tradeDetails")

# Refactor "Renew License": the trade ownership type on this sample row
# changes from "Permanent" to "Temporary".
$ws.Range("C2").Value = "Temporary"

# Update the sheet's last active selection to match.
$ws.Activate()
$ws.Range("G10").Select()
